$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# wdHeaderFooterPrimary = 1, wdHeaderFooterFirstPage = 2
$footerPrimary = $sec.Footers.Item(1)   # -> word/footer2.xml (docPr id="2")
$footerFirst   = $sec.Footers.Item(2)   # -> word/footer1.xml (docPr id="3")
$headerFirst   = $sec.Headers.Item(2)   # -> word/header1.xml (docPr id="1")

# Pearson Edexcel logo in the primary (default) page footer: image2.png -> image1.png
$picFooterPrimary = $footerPrimary.Range.InlineShapes.Item(1)
[void]$picFooterPrimary.Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.png"

# Pearson Edexcel logo in the first-page footer: image2.png -> image1.png
$picFooterFirst = $footerFirst.Range.InlineShapes.Item(1)
[void]$picFooterFirst.Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.png"

# BTec logo in the first-page header: image1.jpg -> image2.jpg
$picHeaderFirst = $headerFirst.Range.InlineShapes.Item(1)
[void]$picHeaderFirst.Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.jpg"

Write-Host "Renamed 3 inline picture(s)."
